# Apply the "MCU-ESP commands" table update:
#  - append the 2-byte end-of-command marker "84 F1" to every command code
#  - drop the LM35 (body temperature) value from the MAX30100 heart-rate /
#    SpO2 command, since it now has its own dedicated command
#  - add a new command row for sending the LM35 body-temperature value,
#    which pushes the old "Relay status" row down to row 13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCU-ESP commands")

$mcuToEsp = "MCU -> ESP"
$espToMcu = "ESP -> MCU"

$rows = @(
    @("84 F0 80 81 84 F1",          $mcuToEsp, "Start wifi configuration mode",                       ""),
    @("84 F0 80 82 84 F1",          $mcuToEsp, "Stop wifi configuration mode",                        ""),
    @("84 F0 80 83 84 F1",          $mcuToEsp, "Reset ESP8266",                                       "Optional"),
    @("84 F0 80 84 84 F1",          $mcuToEsp, "Get Wifi status of ESP8266",                          ""),
    @("84 F0 80 85 84 F1",          $mcuToEsp, "Connect to WiFi",                                     ""),
    @("84 F0 81 82 84 F1",          $espToMcu, "Finish wifi configuration",                           ""),
    @("84 F0 81 84 XX 84 F1",       $espToMcu, "Return the Wifi status of ESP8266",                   ""),
    @("84 F0 82 81 XX YY 84 F1",    $mcuToEsp, "Send DHT11 values: temperature XX and humidity YY",   ""),
    @("84 F0 82 82 XX XX YY 84 F1", $mcuToEsp, "Send Heart rate XXXX, SpO2 YY",                       ""),
    @("84 F0 82 83 XX 84 F1",       $mcuToEsp, "Send the SOS status XX",                               ""),
    @("84 F0 82 84 XX 84 F1",       $mcuToEsp, "Send the body temperature XX",                         ""),
    @("84 F0 83 XX YY 84 F1",       $espToMcu, "Send the status of Relay 1 (XX) and Relay 2 (YY)",     "")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A14").Select() | Out-Null
